$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 15.43697013441807
$ws.Range("F3").Value = 11.09358157116992
$ws.Range("F4").Value = 10.23060859044444
$ws.Range("F5").Value = 10.0348221501642
$ws.Range("F6").Value = 9.495077113351307
$ws.Range("F7").Value = 7.271537858766404

$ws.Range("B8").Value = 7
$ws.Range("C8").Value = "6024c18b094ac71dd93f4f5a"
$ws.Range("D8").Value = "Katherine"
$ws.Range("E8").Value = "female"
$ws.Range("F8").Value = 5.429885939330042

$ws.Range("B9").Value = 9
$ws.Range("C9").Value = "5e35d91ea42bce592e996843"
$ws.Range("D9").Value = "Sergio"
$ws.Range("E9").Value = "male"
$ws.Range("F9").Value = 5.174776003077755

$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "5f0142aa1eb1e528e7abce50"
$ws.Range("D10").Value = "Valeria"
$ws.Range("F10").Value = 5.109387050937287

$ws.Range("F11").Value = 4.066289137297956

$ws.Range("B12").Value = 12
$ws.Range("C12").Value = "5e58b3e415b8d40b5e1dabf1"
$ws.Range("D12").Value = "Cristian"
$ws.Range("F12").Value = 3.287022510508467

$ws.Range("B13").Value = 11
$ws.Range("C13").Value = "5f5ea8227fa75676f56f9276"
$ws.Range("D13").Value = "Carlos"
$ws.Range("F13").Value = 3.121599655266171
